$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -8.281899999999997
$ws.Range("D7").Value = -7.489399999999997
$ws.Range("C8").Value = -11.16599999999999

$ws.Range("A12").Value = -21.76950000000003
$ws.Range("C12").Value = -12.84649999999999

$ws.Range("C14").Value = -12.09689999999999

$ws.Range("D19").Value = -8.590699999999991

$ws.Range("D21").Value = -7.733499999999997

$ws.Range("C22").Value = -10.28799999999999

$ws.Range("D24").Value = -8.69769999999999
